$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style of the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for column H (rows 2-31), 1 for top results, 0 otherwise
$hValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 1
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 0
}

foreach ($row in $hValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $hValues[$row]
}
